$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.075908364421962
$ws.Range("D2").Value = 1.077261357775101
$ws.Range("E2").Value = 1.079356729958614
$ws.Range("F2").Value = 1.089179956340379
$ws.Range("I2").Value = 1.059076599191657
$ws.Range("J2").Value = 1.080810472197302
$ws.Range("K2").Value = 1.079943122504699
$ws.Range("L2").Value = 1.082032994152229
$ws.Range("M2").Value = 1.091830738180567
$ws.Range("N2").Value = 1.082345347038516

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.077243431783943
$ws.Range("D3").Value = 1.078159405063759
$ws.Range("E3").Value = 1.080499628080601
$ws.Range("F3").Value = 1.090289831200039
$ws.Range("I3").Value = 1.059467636986746
$ws.Range("J3").Value = 1.081803366582757
$ws.Range("K3").Value = 1.080658632997403
$ws.Range("L3").Value = 1.082993156911737
$ws.Range("M3").Value = 1.09275980305793
$ws.Range("N3").Value = 1.083339651447886

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.078107060685423
$ws.Range("D4").Value = 1.078740262572961
$ws.Range("E4").Value = 1.081239142488592
$ws.Range("F4").Value = 1.091008005143443
$ws.Range("I4").Value = 1.059719345832128
$ws.Range("J4").Value = 1.082445059564098
$ws.Range("K4").Value = 1.081120738932602
$ws.Range("L4").Value = 1.083613839765412
$ws.Range("M4").Value = 1.093360382327046
$ws.Range("N4").Value = 1.083982255706863

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.07847007353253
$ws.Range("D5").Value = 1.078984398462305
$ws.Range("E5").Value = 1.081550031848547
$ws.Range("F5").Value = 1.091309929295947
$ws.Range("I5").Value = 1.059824849150464
$ws.Range("J5").Value = 1.082714643292081
$ws.Range("K5").Value = 1.081314799060352
$ws.Range("L5").Value = 1.083874631111688
$ws.Range("M5").Value = 1.093612726356936
$ws.Range("N5").Value = 1.084252222274661

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.078531021742498
$ws.Range("D6").Value = 1.079025386651493
$ws.Range("E6").Value = 1.081602231440005
$ws.Range("F6").Value = 1.091360623950854
$ws.Range("I6").Value = 1.059842545163776
$ws.Range("J6").Value = 1.082759896848412
$ws.Range("K6").Value = 1.081347370364772
$ws.Range("L6").Value = 1.083918410741976
$ws.Range("M6").Value = 1.093655087885684
$ws.Range("N6").Value = 1.084297540096232

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.078111911500241
$ws.Range("D7").Value = 1.078743524950295
$ws.Range("E7").Value = 1.081243296615254
$ws.Range("F7").Value = 1.091012039452129
$ws.Range("I7").Value = 1.059720756809891
$ws.Range("J7").Value = 1.082448662475828
$ws.Range("K7").Value = 1.08112333279458
$ws.Range("L7").Value = 1.083617325035007
$ws.Range("M7").Value = 1.093363754706906
$ws.Range("N7").Value = 1.083985863735141

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.076359609955664
$ws.Range("D8").Value = 1.077564906964118
$ws.Range("E8").Value = 1.079742982232745
$ws.Range("F8").Value = 1.089555042294873
$ws.Range("I8").Value = 1.059209025683535
$ws.Range("J8").Value = 1.081146187276087
$ws.Range("K8").Value = 1.080185114735373
$ws.Range("L8").Value = 1.082357612173658
$ws.Range("M8").Value = 1.092144842463351
$ws.Range("N8").Value = 1.08268153887122

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.073269796003791
$ws.Range("D9").Value = 1.075486170599834
$ws.Range("E9").Value = 1.077099033738972
$ws.Range("F9").Value = 1.086987651845554
$ws.Range("I9").Value = 1.058297161427716
$ws.Range("J9").Value = 1.078845035135648
$ws.Range("K9").Value = 1.078525100288056
$ws.Range("L9").Value = 1.080133122982036
$ws.Range("M9").Value = 1.089992405631533
$ws.Range("N9").Value = 1.080377118830796

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.071208371041656
$ws.Range("D10").Value = 1.074099051377085
$ws.Range("E10").Value = 1.075336162034588
$ws.Range("F10").Value = 1.085275992149388
$ws.Range("I10").Value = 1.057682400828378
$ws.Range("J10").Value = 1.077306768494726
$ws.Range("K10").Value = 1.077413830157047
$ws.Range("L10").Value = 1.078646863511584
$ws.Range("M10").Value = 1.08855430326012
$ws.Range("N10").Value = 1.078836667674802

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.070315343002008
$ws.Range("D11").Value = 1.073498095675695
$ws.Range("E11").Value = 1.074572738660173
$ws.Range("F11").Value = 1.084534789699484
$ws.Range("I11").Value = 1.05741456813263
$ws.Range("J11").Value = 1.076639667940609
$ws.Range("K11").Value = 1.076931535938136
$ws.Range("L11").Value = 1.078002499643693
$ws.Range("M11").Value = 1.087930824658658
$ws.Range("N11").Value = 1.078168619761381

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.069983566302618
$ws.Range("D12").Value = 1.073274824290128
$ws.Range("E12").Value = 1.074289153826329
$ws.Range("F12").Value = 1.084259465964029
$ws.Range("I12").Value = 1.05731483626787
$ws.Range("J12").Value = 1.076391721358779
$ws.Range("K12").Value = 1.076752222931966
$ws.Range("L12").Value = 1.07776303145219
$ws.Range("M12").Value = 1.087699119448968
$ws.Range("N12").Value = 1.077920321066965

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.070054736590207
$ws.Range("D13").Value = 1.073322719024764
$ws.Range("E13").Value = 1.074349984467158
$ws.Range("F13").Value = 1.084318524236708
$ws.Range("I13").Value = 1.057336240277406
$ws.Range("J13").Value = 1.076444913827288
$ws.Range("K13").Value = 1.076790693776998
$ws.Range("L13").Value = 1.077814403779726
$ws.Range("M13").Value = 1.087748826347252
$ws.Range("N13").Value = 1.07797358907488

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.070287919610423
$ws.Range("D14").Value = 1.073479641012675
$ws.Range("E14").Value = 1.074549297775904
$ws.Range("F14").Value = 1.084512031529941
$ws.Range("I14").Value = 1.057406329306938
$ws.Range("J14").Value = 1.076619175789324
$ws.Range("K14").Value = 1.076916717291935
$ws.Range("L14").Value = 1.077982707644224
$ws.Range("M14").Value = 1.087911674240468
$ws.Range("N14").Value = 1.07814809850889

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.070431582468575
$ws.Range("D15").Value = 1.073576319214884
$ws.Range("E15").Value = 1.074672099178759
$ws.Range("F15").Value = 1.084631256628052
$ws.Range("I15").Value = 1.057449480731079
$ws.Range("J15").Value = 1.076726523628982
$ws.Range("K15").Value = 1.076994342314128
$ws.Range("L15").Value = 1.07808638888258
$ws.Range("M15").Value = 1.088011994581866
$ws.Range("N15").Value = 1.078255598794794

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.071267629268905
$ws.Range("D16").Value = 1.074138927916642
$ws.Range("E16").Value = 1.075386825887285
$ws.Range("F16").Value = 1.085325182269506
$ws.Range("I16").Value = 1.057700141430325
$ws.Range("J16").Value = 1.077351020051046
$ws.Range("K16").Value = 1.07744581501581
$ws.Range("L16").Value = 1.078689610714481
$ws.Range("M16").Value = 1.088595665097444
$ws.Range("N16").Value = 1.078880982073407

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.071791945035184
$ws.Range("D17").Value = 1.07449175005311
$ws.Range("E17").Value = 1.075835129923892
$ws.Range("F17").Value = 1.085760450872085
$ws.Range("I17").Value = 1.057856935284717
$ws.Range("J17").Value = 1.077742475165235
$ws.Range("K17").Value = 1.077728714744596
$ws.Range("L17").Value = 1.079067779237725
$ws.Range("M17").Value = 1.088961578635439
$ws.Range("N17").Value = 1.079272993098759

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.072097729520561
$ws.Range("D18").Value = 1.074697514052067
$ws.Range("E18").Value = 1.076096609710564
$ws.Range("F18").Value = 1.08601433170288
$ws.Range("I18").Value = 1.057948232572411
$ws.Range("J18").Value = 1.077970705963722
$ws.Range("K18").Value = 1.077893618647834
$ws.Range("L18").Value = 1.079288281130211
$ws.Range("M18").Value = 1.089174935557724
$ws.Range("N18").Value = 1.079501548011158

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.072201987446659
$ws.Range("D19").Value = 1.074767668999096
$ws.Range("E19").Value = 1.076185766204382
$ws.Range("F19").Value = 1.086100897882652
$ws.Range("I19").Value = 1.057979335829773
$ws.Range("J19").Value = 1.078048510157192
$ws.Range("K19").Value = 1.077949828549953
$ws.Range("L19").Value = 1.07936345350773
$ws.Range("M19").Value = 1.089247672210641
$ws.Range("N19").Value = 1.079579462695508

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.071735695083554
$ws.Range("D20").Value = 1.074453898810608
$ws.Range("E20").Value = 1.075787032026088
$ws.Range("F20").Value = 1.085713751095876
$ws.Range("I20").Value = 1.057840129128838
$ws.Range("J20").Value = 1.077700485926675
$ws.Range("K20").Value = 1.077698373328934
$ws.Range("L20").Value = 1.079027213356454
$ws.Range("M20").Value = 1.088922327269383
$ws.Range("N20").Value = 1.079230944230665

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.070219254882785
$ws.Range("D21").Value = 1.073433432778022
$ws.Range("E21").Value = 1.074490605416573
$ws.Range("F21").Value = 1.084455048694079
$ws.Range("I21").Value = 1.057385696656821
$ws.Range("J21").Value = 1.076567864302759
$ws.Range("K21").Value = 1.076879611139849
$ws.Range("L21").Value = 1.077933149755982
$ws.Range("M21").Value = 1.087863722851923
$ws.Range("N21").Value = 1.078096714154129

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.069265423221056
$ws.Range("D22").Value = 1.072791537337385
$ws.Range("E22").Value = 1.073675400382524
$ws.Range("F22").Value = 1.083663604544688
$ws.Range("I22").Value = 1.057098548132589
$ws.Range("J22").Value = 1.07585483794163
$ws.Range("K22").Value = 1.076363853411399
$ws.Range("L22").Value = 1.07724455729054
$ws.Range("M22").Value = 1.087197455594955
$ws.Range("N22").Value = 1.077382675213783

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.069771105181404
$ws.Range("D23").Value = 1.073131845889301
$ws.Range("E23").Value = 1.074107565256758
$ws.Range("F23").Value = 1.084083169199988
$ws.Range("I23").Value = 1.057250906728913
$ws.Range("J23").Value = 1.076232913008848
$ws.Range("K23").Value = 1.076637358620951
$ws.Range("L23").Value = 1.077609661276927
$ws.Range("M23").Value = 1.087550721558235
$ws.Range("N23").Value = 1.077761287190962

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.071761112139044
$ws.Range("D24").Value = 1.074471002254274
$ws.Range("E24").Value = 1.075808765421664
$ws.Range("F24").Value = 1.085734852726951
$ws.Range("I24").Value = 1.05784772359522
$ws.Range("J24").Value = 1.077719459360035
$ws.Range("K24").Value = 1.077712083638457
$ws.Range("L24").Value = 1.079045543571258
$ws.Range("M24").Value = 1.088940063501626
$ws.Range("N24").Value = 1.079249944608476

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.074068846710272
$ws.Range("D25").Value = 1.076023798788963
$ws.Range("E25").Value = 1.077782592420476
$ws.Range("F25").Value = 1.087651389377916
$ws.Range("I25").Value = 1.058534104722807
$ws.Range("J25").Value = 1.079440663841546
$ws.Range("K25").Value = 1.078955059635327
$ws.Range("L25").Value = 1.08070877539442
$ws.Range("M25").Value = 1.090549410489585
$ws.Range("N25").Value = 1.080973593397776
